$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-359) holds the "Förändrad" (last changed) date serial,
# which was bumped by one day (45178 -> 45179) for every data row.
$ws.Range("C2:C359").Value = 45179
